# Helper: write a value as genuine text (no number auto-detection), leaving
# no residual style/numberformat on the cell.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Item("2022-Q1")

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right before the existing "2022-Q1"
#    sheet, so the tab order becomes: 总计, 2022-Q3, 2022-Q1.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($q1)
$q3.Name = "2022-Q3"

# ------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: push the existing 2022-Q1 row
#    down and insert a new row for 2022-Q3 above it.
# ------------------------------------------------------------------
$total.Rows(2).Insert()

Set-TextValue $total.Cells.Item(2,2) "2022-Q3"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,3).ClearFormats()
$total.Cells.Item(2,4).Value = 1.38
$total.Cells.Item(2,4).ClearFormats()
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(1,2).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

# Row 3 (previously row 2, shifted down by the insert) keeps its data but
# its index value moves from 0 to 1.
$total.Cells.Item(3,1).Value = 1

# ------------------------------------------------------------------
# 3. Populate the new "2022-Q3" sheet with the fund-holdings detail.
# ------------------------------------------------------------------
Set-TextValue $q3.Cells.Item(1,2) "基金代码"
Set-TextValue $q3.Cells.Item(1,3) "基金名称"
Set-TextValue $q3.Cells.Item(1,4) "基金规模"
Set-TextValue $q3.Cells.Item(1,5) "股票总仓位"
Set-TextValue $q3.Cells.Item(1,6) "仓位占比"
Set-TextValue $q3.Cells.Item(1,7) "持有市值(亿元)"
Set-TextValue $q3.Cells.Item(1,8) "仓位排名"

$rows = @(
    @(0, "002350", "华安安华灵活配置混合A", "32.74", "87.55", "1.96", "0.6417", 10),
    @(1, "160212", "国泰估值优势混合（LOF）A", "9.14", "94.29", "4.84", "0.4424", 8),
    @(2, "007731", "民生加银持续成长混合A", "3.22", "94.57", "5.68", "0.1829", 8),
    @(3, "007732", "民生加银持续成长混合C", "1.89", "94.57", "5.68", "0.1074", 8),
    @(4, "016183", "华安安华灵活配置混合C", "0.45", "87.55", "1.96", "0.0088", 10),
    @(5, "016616", "国泰估值优势混合（LOF）C", "0.00", "94.29", "4.84", $null, 8)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r,1).Value = $row[0]
    Set-TextValue $q3.Cells.Item($r,2) $row[1]
    Set-TextValue $q3.Cells.Item($r,3) $row[2]
    Set-TextValue $q3.Cells.Item($r,4) $row[3]
    Set-TextValue $q3.Cells.Item($r,5) $row[4]
    Set-TextValue $q3.Cells.Item($r,6) $row[5]
    if ($row[6] -eq $null) {
        $q3.Cells.Item($r,7).Value = 0
    } else {
        Set-TextValue $q3.Cells.Item($r,7) $row[6]
    }
    $q3.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# Reapply the bold/bordered "index" style (matching style used on the
# 总计 sheet) to the header row and the A-column index column.
$total.Cells.Item(1,2).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Cells.Item(1,2).Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)

# The "2022-Q1" sheet (now the 3rd tab) was, and remains, the selected /
# active sheet - restore that (creating the new sheet shifts focus away).
# Re-fetch by name: the sheet was re-indexed/re-positioned by the Add()
# call above, so the original $q1 reference is stale.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Activate()
